# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# ----- Update "last updated" timestamp text (A1) -----
$ws.Range("A1").Value = "Datos actualizados a 28 de Abril de 2020 a las 12:22"

# ----- Reorder country names (text swaps caused by new ranking) -----
# Catar / Emiratos Arabes Unidos swap
$ws.Range("A36").Value = "Emiratos Arabes Unidos"
$ws.Range("A37").Value = "Catar"

# Albania / Principado de Andorra / Senegal / Libano / Crucero reshuffle
$ws.Range("A94").Value = "Albania"
$ws.Range("A95").Value = "Principado de Andorra"
$ws.Range("A96").Value = "Senegal"
$ws.Range("A97").Value = "Libano"
$ws.Range("A98").Value = "Crucero"

# ----- Updated case numbers -----
$ws.Range("E19").Value = 5287
$ws.Range("G19").Value = 12
$ws.Range("H19").Value = 1677

$ws.Range("B35").Value = 11616
$ws.Range("C35").Value = 277
$ws.Range("D35").Value = 3404
$ws.Range("E35").Value = 7562
$ws.Range("F35").Value = 243

$ws.Range("B36").Value = 11380
$ws.Range("C36").Value = 541
$ws.Range("D36").Value = 2181
$ws.Range("E36").Value = 9110
$ws.Range("F36").Value = 1
$ws.Range("G36").Value = 7
$ws.Range("H36").Value = 89

$ws.Range("B37").Value = 11244
$ws.Range("D37").Value = 1066
$ws.Range("E37").Value = 10168
$ws.Range("F37").Value = 72
$ws.Range("H37").Value = 10

$ws.Range("B55").Value = 4246
$ws.Range("C55").Value = 126
$ws.Range("D55").Value = 739
$ws.Range("E55").Value = 3344
$ws.Range("G55").Value = 1
$ws.Range("H55").Value = 163

$ws.Range("B76").Value = 1585
$ws.Range("C76").Value = 20
$ws.Range("D76").Value = 682
$ws.Range("E76").Value = 840
$ws.Range("G76").Value = 3
$ws.Range("H76").Value = 63

$ws.Range("F82").Value = 39

$ws.Range("D88").Value = 811
$ws.Range("E88").Value = 223

$ws.Range("B94").Value = 750
$ws.Range("C94").Value = 14
$ws.Range("D94").Value = 431
$ws.Range("E94").Value = 289
$ws.Range("F94").Value = 4
$ws.Range("G94").Value = 2
$ws.Range("H94").Value = 30

$ws.Range("B95").Value = 743
$ws.Range("D95").Value = 385
$ws.Range("E95").Value = 318
$ws.Range("F95").Value = 17
$ws.Range("H95").Value = 40

$ws.Range("D96").Value = 284
$ws.Range("E96").Value = 443
$ws.Range("F96").Value = 1
$ws.Range("H96").Value = 9

$ws.Range("B97").Value = 717
$ws.Range("C97").Value = 7
$ws.Range("D97").Value = 145
$ws.Range("E97").Value = 548
$ws.Range("F97").Value = 44
$ws.Range("H97").Value = 24

$ws.Range("B98").Value = 712
$ws.Range("D98").Value = 645
$ws.Range("E98").Value = 54
$ws.Range("F98").Value = 4
$ws.Range("H98").Value = 13

$ws.Range("B105").Value = 592
$ws.Range("C105").Value = 4
$ws.Range("E105").Value = 451
